# Applies the per-cell numeric updates recorded for Sheets/Hyperion_Profits.xlsx.
# Workbook has 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); each sheet is a
# "Leve profit" table with columns H..N holding computed market-price/profit figures
# that were recalculated upstream (market prices refreshed) and re-synced here.
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3908.4412
$ws.Range("J17").Value = 4204.0967
$ws.Range("L17").Value = 12612.2901
$ws.Range("N17").Value = -12948.2901
# Row 19
$ws.Range("H19").Value = 2789.8635
$ws.Range("I19").Value = 1018
$ws.Range("K19").Value = 1018
$ws.Range("M19").Value = -843
# Row 76
$ws.Range("H76").Value = 4434.3887
$ws.Range("J76").Value = 4999
$ws.Range("L76").Value = 4999
$ws.Range("N76").Value = -5629
# Row 79
$ws.Range("H79").Value = 4434.3887
$ws.Range("J79").Value = 4999
$ws.Range("L79").Value = 4999
$ws.Range("N79").Value = -7183
# Row 132
$ws.Range("H132").Value = 23812438
$ws.Range("I132").Value = 25002996
$ws.Range("J132").Value = 1299.5
$ws.Range("K132").Value = 75008988
$ws.Range("L132").Value = 3898.5
$ws.Range("M132").Value = -75006458
$ws.Range("N132").Value = -8958.5
# Row 137
$ws.Range("H137").Value = 131108.92
$ws.Range("I137").Value = 359555
$ws.Range("J137").Value = 4194.4443
$ws.Range("K137").Value = 1078665
$ws.Range("L137").Value = 12583.3329
$ws.Range("M137").Value = -1076115
$ws.Range("N137").Value = -17683.3329
# Row 138
$ws.Range("H138").Value = 2632.092
$ws.Range("I138").Value = 1633.7142
$ws.Range("J138").Value = 2857.5322
$ws.Range("K138").Value = 4901.142599999999
$ws.Range("L138").Value = 8572.596600000001
$ws.Range("M138").Value = 238.8574000000008
$ws.Range("N138").Value = -18852.5966
# Row 141
$ws.Range("H141").Value = 2481.6
$ws.Range("I141").Value = 2344.3333
$ws.Range("K141").Value = 7032.999899999999
$ws.Range("M141").Value = -1852.999899999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 1215.1428
$ws.Range("I16").Value = 377
$ws.Range("K16").Value = 377
$ws.Range("M16").Value = -90
# Row 32
$ws.Range("H32").Value = 4560.6724
$ws.Range("I32").Value = 3803.5686
$ws.Range("K32").Value = 3803.5686
$ws.Range("M32").Value = -3516.5686
# Row 45
$ws.Range("H45").Value = 6257643.5
$ws.Range("I45").Value = 11067301
$ws.Range("K45").Value = 11067301
$ws.Range("M45").Value = -11066924
# Row 61
$ws.Range("H61").Value = 2105.2727
$ws.Range("I61").Value = 1816.3
$ws.Range("K61").Value = 1816.3
$ws.Range("M61").Value = -1604.3
# Row 88
$ws.Range("H88").Value = 5005.25
$ws.Range("I88").Value = 2503
$ws.Range("K88").Value = 2503
$ws.Range("M88").Value = -2097
# Row 91
$ws.Range("H91").Value = 5005.25
$ws.Range("I91").Value = 2503
$ws.Range("K91").Value = 2503
$ws.Range("M91").Value = -1099
# Row 122
$ws.Range("H122").Value = 522832.56
$ws.Range("I122").Value = 1633.9706
$ws.Range("K122").Value = 4901.9118
$ws.Range("M122").Value = -2451.9118
# Row 132
$ws.Range("H132").Value = 2616.4075
$ws.Range("I132").Value = 1940.0625
$ws.Range("K132").Value = 5820.1875
$ws.Range("M132").Value = -3290.1875
# Row 136
$ws.Range("H136").Value = 2105.2727
$ws.Range("I136").Value = 1816.3
$ws.Range("K136").Value = 5448.9
$ws.Range("M136").Value = -2898.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6676436
$ws.Range("I86").Value = 6676436
$ws.Range("K86").Value = 6676436
$ws.Range("M86").Value = -6675313
# Row 89
$ws.Range("H89").Value = 6676436
$ws.Range("I89").Value = 6676436
$ws.Range("K89").Value = 33382180
$ws.Range("M89").Value = -33376564
# Row 134
$ws.Range("H134").Value = 2496.0754
$ws.Range("I134").Value = 1034.1464
$ws.Range("K134").Value = 3102.4392
$ws.Range("M134").Value = -567.4392000000003

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 31283.742
$ws.Range("I31").Value = 1476.0952
$ws.Range("J31").Value = 93879.8
$ws.Range("K31").Value = 1476.0952
$ws.Range("L31").Value = 93879.8
$ws.Range("M31").Value = -1181.0952
$ws.Range("N31").Value = -94469.8
# Row 34
$ws.Range("H34").Value = 31283.742
$ws.Range("I34").Value = 1476.0952
$ws.Range("J34").Value = 93879.8
$ws.Range("K34").Value = 1476.0952
$ws.Range("L34").Value = 93879.8
$ws.Range("M34").Value = -1274.0952
$ws.Range("N34").Value = -94283.8
# Row 35
$ws.Range("H35").Value = 5309.875
$ws.Range("I35").Value = 1096
$ws.Range("K35").Value = 1096
$ws.Range("M35").Value = -802
# Row 62
$ws.Range("H62").Value = 1874.25
$ws.Range("I62").Value = 1499
$ws.Range("K62").Value = 1499
$ws.Range("M62").Value = -875
# Row 65
$ws.Range("H65").Value = 1874.25
$ws.Range("I65").Value = 1499
$ws.Range("K65").Value = 7495
$ws.Range("M65").Value = -4375
# Row 105
$ws.Range("H105").Value = 1692.3846
$ws.Range("I105").Value = 1680.25
$ws.Range("J105").Value = 1732.8334
$ws.Range("K105").Value = 1680.25
$ws.Range("L105").Value = 1732.8334
$ws.Range("M105").Value = 66.75
$ws.Range("N105").Value = -5226.8334
# Row 132
$ws.Range("H132").Value = 80330.96000000001
$ws.Range("I132").Value = 49160
$ws.Range("K132").Value = 147480
$ws.Range("M132").Value = -144950
# Row 134
$ws.Range("H134").Value = 35723.68
$ws.Range("I134").Value = 71727.234
$ws.Range("J134").Value = 4520.6
$ws.Range("K134").Value = 215181.702
$ws.Range("L134").Value = 13561.8
$ws.Range("M134").Value = -212646.702
$ws.Range("N134").Value = -18631.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Range("H50").Value = 837.8889
$ws.Range("I50").Value = 236.2
$ws.Range("J50").Value = 1590
$ws.Range("K50").Value = 708.5999999999999
$ws.Range("L50").Value = 4770
$ws.Range("M50").Value = -227.5999999999999
$ws.Range("N50").Value = -5732
# Row 53
$ws.Range("H53").Value = 837.8889
$ws.Range("I53").Value = 236.2
$ws.Range("J53").Value = 1590
$ws.Range("K53").Value = 708.5999999999999
$ws.Range("L53").Value = 4770
$ws.Range("M53").Value = -227.5999999999999
$ws.Range("N53").Value = -5732

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 6799840
$ws.Range("I102").Value = 9262691
$ws.Range("K102").Value = 9262691
$ws.Range("M102").Value = -9261069
# Row 113
$ws.Range("H113").Value = 20834966
$ws.Range("I113").Value = 23811090
$ws.Range("K113").Value = 23811090
$ws.Range("M113").Value = -23808920
# Row 126
$ws.Range("H126").Value = 4956484.5
$ws.Range("I126").Value = 3249374.2
$ws.Range("K126").Value = 9748122.600000001
$ws.Range("M126").Value = -9745652.600000001
# Row 132
$ws.Range("H132").Value = 3723.4707
$ws.Range("I132").Value = 2957.6924
$ws.Range("J132").Value = 6212.25
$ws.Range("K132").Value = 8873.0772
$ws.Range("L132").Value = 18636.75
$ws.Range("M132").Value = -6343.0772
$ws.Range("N132").Value = -23696.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 822.0526
$ws.Range("I16").Value = 859.3889
$ws.Range("K16").Value = 859.3889
$ws.Range("M16").Value = -689.3889
# Row 46
$ws.Range("H46").Value = 4149.8335
$ws.Range("I46").Value = 2225
$ws.Range("K46").Value = 2225
$ws.Range("M46").Value = -2037
# Row 55
$ws.Range("H55").Value = 1626.4242
$ws.Range("I55").Value = 1265.421
$ws.Range("J55").Value = 2116.3572
$ws.Range("K55").Value = 1265.421
$ws.Range("L55").Value = 2116.3572
$ws.Range("M55").Value = -1092.421
$ws.Range("N55").Value = -2462.3572
# Row 61
$ws.Range("H61").Value = 3270710.5
$ws.Range("I61").Value = 5294091
$ws.Range("J61").Value = 2173
$ws.Range("K61").Value = 5294091
$ws.Range("L61").Value = 2173
$ws.Range("M61").Value = -5293889
$ws.Range("N61").Value = -2577
# Row 82
$ws.Range("H82").Value = 3969513.5
$ws.Range("I82").Value = 6945512.5
$ws.Range("J82").Value = 1515
$ws.Range("K82").Value = 6945512.5
$ws.Range("L82").Value = 1515
$ws.Range("M82").Value = -6945151.5
$ws.Range("N82").Value = -2237
# Row 85
$ws.Range("H85").Value = 3969513.5
$ws.Range("I85").Value = 6945512.5
$ws.Range("J85").Value = 1515
$ws.Range("K85").Value = 6945512.5
$ws.Range("L85").Value = 1515
$ws.Range("M85").Value = -6944264.5
$ws.Range("N85").Value = -4011
# Row 113
$ws.Range("H113").Value = 3270710.5
$ws.Range("I113").Value = 5294091
$ws.Range("J113").Value = 2173
$ws.Range("K113").Value = 5294091
$ws.Range("L113").Value = 2173
$ws.Range("M113").Value = -5291921
$ws.Range("N113").Value = -6513
# Row 122
$ws.Range("H122").Value = 4716.3667
$ws.Range("I122").Value = 3194.625
$ws.Range("J122").Value = 6455.5
$ws.Range("K122").Value = 9583.875
$ws.Range("L122").Value = 19366.5
$ws.Range("M122").Value = -7133.875
$ws.Range("N122").Value = -24266.5
# Row 132
$ws.Range("H132").Value = 4811.7925
$ws.Range("I132").Value = 4115.6763
$ws.Range("K132").Value = 12347.0289
$ws.Range("M132").Value = -9817.028900000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1297.5883
$ws.Range("I100").Value = 543.8461
$ws.Range("K100").Value = 1087.6922
$ws.Range("M100").Value = -546.6922
# Row 122
$ws.Range("H122").Value = 2428.6667
$ws.Range("I122").Value = 865.6667
$ws.Range("K122").Value = 2597.0001
$ws.Range("M122").Value = -147.0001000000002
# Row 123
$ws.Range("H123").Value = 59966.332
$ws.Range("J123").Value = 59966.332
$ws.Range("L123").Value = 59966.332
$ws.Range("N123").Value = -69766.33199999999
# Row 126
$ws.Range("H126").Value = 1654.619
$ws.Range("I126").Value = 1955.2727
$ws.Range("J126").Value = 1323.9
$ws.Range("K126").Value = 5865.8181
$ws.Range("L126").Value = 3971.7
$ws.Range("M126").Value = -3395.8181
$ws.Range("N126").Value = -8911.700000000001
# Row 132
$ws.Range("H132").Value = 50547320
$ws.Range("I132").Value = 66668860
$ws.Range("K132").Value = 200006580
$ws.Range("M132").Value = -200004050
# Row 136
$ws.Range("H136").Value = 5137.125
$ws.Range("I136").Value = 3400
$ws.Range("K136").Value = 10200
$ws.Range("M136").Value = -7650

Write-Output "Applied 247 cell updates across 8 sheets."
